# "Loan RBI, Variable Instalments"
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted right
# before column N, pushing the old N:P columns (Late / heading / Outstanding)
# one position to the right, into O:Q. The freshly inserted column N is left
# blank (to later hold per-instalment variable amounts), and is given an
# explicit width. The edit finishes with "Repayment schedule" as the active
# sheet/tab (it was previously "Edit Repayment Schedule"), and updated cell
# selections on both of those sheets.

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsEdit  = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new blank column before column N (existing N/O/P shift right to O/P/Q).
$wsRepay.Activate()
$wsRepay.Columns("N").Insert() | Out-Null
$wsRepay.Columns("N").ColumnWidth = 10.14

# The previously-active sheet ("Edit Repayment Schedule") now has a different
# selected cell too.
$wsEdit.Activate()
$wsEdit.Range("E16").Select() | Out-Null

# Leave "Repayment schedule" as the active sheet/tab with its new selection.
$wsRepay.Activate()
$wsRepay.Range("L15").Select() | Out-Null
